$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking values (mirrors the source data
# being written as literal strings) by applying a Text number format
# before assigning each value, then writing the exact literal value.
$cellValues = [ordered]@{
    "D2" = "305.64"
    "E2" = "-4.92%"
    "G2" = "6"
    "D3" = "39.62"
    "E3" = "-7.85%"
    "G3" = "6"
    "D4" = "5.037"
    "E4" = "-2.84%"
    "G4" = "6"
    "D5" = "0.07677"
    "E5" = "-5.88%"
    "G5" = "6"
    "D6" = "4.266"
    "E6" = "-1.60%"
    "G6" = "6"
    "D7" = "1.616"
    "E7" = "-10.48%"
    "G7" = "6"
    "D8" = "0.8818"
    "E8" = "-7.02%"
    "G8" = "6"
    "D9" = "0.09651"
    "E9" = "-13.43%"
    "G9" = "6"
    "D10" = "0.1717"
    "E10" = "-7.66%"
    "G10" = "6"
    "D11" = "0.04478"
    "E11" = "-3.40%"
    "G11" = "6"
    "D12" = "0.08860"
    "E12" = "-5.57%"
    "G12" = "6"
    "D13" = "0.1055"
    "E13" = "-0.42%"
    "G13" = "6"
    "D14" = "0.001267"
    "E14" = "-1.69%"
    "G14" = "6"
    "D15" = "0.005945"
    "E15" = "5.31%"
    "G15" = "6"
    "D16" = "3.352"
    "E16" = "-0.34%"
    "G16" = "6"
    "D17" = "2.452"
    "E17" = "-2.47%"
    "G17" = "6"
    "D18" = "0.3300"
    "E18" = "-1.94%"
    "G18" = "6"
    "D19" = "6.978"
    "E19" = "-5.88%"
    "G19" = "6"
    "D20" = "0.1349"
    "E20" = "-2.93%"
    "G20" = "6"
    "D21" = "0.3224"
    "E21" = "22.69%"
    "G21" = "6"
    "D22" = "0.04207"
    "E22" = "0.90%"
    "G22" = "6"
    "D23" = "0.001194"
    "E23" = "-4.50%"
    "G23" = "6"
    "D24" = "0.004062"
    "E24" = "-5.45%"
    "G24" = "6"
    "D25" = "0.0001221"
    "G25" = "6"
    "E26" = "-0.09%"
    "G26" = "6"
    "G27" = "6"
    "G28" = "6"
    "G29" = "6"
    "G30" = "6"
    "G31" = "6"
    "G32" = "6"
    "G33" = "6"
    "G34" = "6"
    "G35" = "6"
    "G36" = "6"
    "G37" = "6"
    "D38" = "0.02325"
    "E38" = "-13.05%"
    "G38" = "6"
    "D39" = "0.05143"
    "E39" = "-6.98%"
    "G39" = "6"
    "D40" = "0.007917"
    "E40" = "-0.40%"
    "G40" = "6"
    "E41" = "-4.86%"
    "G41" = "6"
    "D42" = "0.006347"
    "E42" = "-3.22%"
    "G42" = "6"
    "D43" = "0.001931"
    "E43" = "-8.93%"
    "G43" = "6"
    "D44" = "0.008694"
    "E44" = "17.82%"
    "G44" = "6"
    "D45" = "0.3031"
    "E45" = "-5.50%"
    "G45" = "6"
    "D46" = "0.00006509"
    "E46" = "-7.00%"
    "G46" = "6"
    "D47" = "0.00000000750"
    "E47" = "0.02%"
    "G47" = "6"
    "D48" = "0.007003"
    "E48" = "98.02%"
    "G48" = "6"
    "D49" = "0.003370"
    "E49" = "-3.06%"
    "G49" = "6"
    "D50" = "0.00002101"
    "E50" = "0.02%"
    "G50" = "6"
    "D51" = "0.0002001"
    "E51" = "0.02%"
    "G51" = "6"
}

foreach ($ref in $cellValues.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $cellValues[$ref]
}
